# Auto-generated: apply numeric cell updates per the commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 3817.92
$ws.Cells.Item(15, 9).Value = 3817.92
$ws.Cells.Item(15, 11).Value = 11453.76
$ws.Cells.Item(15, 13).Value = -11284.76
$ws.Cells.Item(53, 8).Value = 1337.9375
$ws.Cells.Item(53, 9).Value = 2247.25
$ws.Cells.Item(53, 11).Value = 2247.25
$ws.Cells.Item(53, 13).Value = -1610.25
$ws.Cells.Item(76, 8).Value = 7500
$ws.Cells.Item(76, 9).Value = 7500
$ws.Cells.Item(76, 11).Value = 7500
$ws.Cells.Item(76, 13).Value = -7185
$ws.Cells.Item(79, 8).Value = 7500
$ws.Cells.Item(79, 9).Value = 7500
$ws.Cells.Item(79, 11).Value = 7500
$ws.Cells.Item(79, 13).Value = -6408
$ws.Cells.Item(86, 8).Value = 1998
$ws.Cells.Item(86, 9).Value = 1998
$ws.Cells.Item(86, 11).Value = 1998
$ws.Cells.Item(86, 13).Value = -875
$ws.Cells.Item(89, 8).Value = 1998
$ws.Cells.Item(89, 9).Value = 1998
$ws.Cells.Item(89, 11).Value = 9990
$ws.Cells.Item(89, 13).Value = -4374
$ws.Cells.Item(105, 8).Value = 60000
$ws.Cells.Item(105, 10).Value = 60000
$ws.Cells.Item(105, 12).Value = 60000
$ws.Cells.Item(105, 14).Value = -66988
$ws.Cells.Item(137, 8).Value = 1627.1
$ws.Cells.Item(137, 10).Value = 1885.4445
$ws.Cells.Item(137, 12).Value = 5656.333500000001
$ws.Cells.Item(137, 14).Value = -10756.3335

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3776.56
$ws.Cells.Item(45, 9).Value = 2744.0625
$ws.Cells.Item(45, 10).Value = 5612.1113
$ws.Cells.Item(45, 11).Value = 2744.0625
$ws.Cells.Item(45, 12).Value = 5612.1113
$ws.Cells.Item(45, 13).Value = -2367.0625
$ws.Cells.Item(45, 14).Value = -6366.1113
$ws.Cells.Item(61, 8).Value = 3150741.8
$ws.Cells.Item(61, 9).Value = 4390831
$ws.Cells.Item(61, 11).Value = 4390831
$ws.Cells.Item(61, 13).Value = -4390619
$ws.Cells.Item(74, 8).Value = 4236.3667
$ws.Cells.Item(74, 9).Value = 2279.6
$ws.Cells.Item(74, 11).Value = 2279.6
$ws.Cells.Item(74, 13).Value = -1405.6
$ws.Cells.Item(77, 8).Value = 4236.3667
$ws.Cells.Item(77, 9).Value = 2279.6
$ws.Cells.Item(77, 11).Value = 11398
$ws.Cells.Item(77, 13).Value = -7030
$ws.Cells.Item(132, 8).Value = 3010.1084
$ws.Cells.Item(132, 9).Value = 2347.6572
$ws.Cells.Item(132, 11).Value = 7042.971600000001
$ws.Cells.Item(132, 13).Value = -4512.971600000001
$ws.Cells.Item(136, 8).Value = 3150741.8
$ws.Cells.Item(136, 9).Value = 4390831
$ws.Cells.Item(136, 11).Value = 13172493
$ws.Cells.Item(136, 13).Value = -13169943

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(70, 8).Value = 150000
$ws.Cells.Item(70, 10).Value = 150000
$ws.Cells.Item(70, 12).Value = 150000
$ws.Cells.Item(70, 14).Value = -150586
$ws.Cells.Item(73, 8).Value = 150000
$ws.Cells.Item(73, 10).Value = 150000
$ws.Cells.Item(73, 12).Value = 150000
$ws.Cells.Item(73, 14).Value = -152028
$ws.Cells.Item(80, 8).Value = 528.0714
$ws.Cells.Item(80, 9).Value = 269.14285
$ws.Cells.Item(80, 10).Value = 787
$ws.Cells.Item(80, 11).Value = 269.14285
$ws.Cells.Item(80, 12).Value = 787
$ws.Cells.Item(80, 13).Value = 728.85715
$ws.Cells.Item(80, 14).Value = -2783
$ws.Cells.Item(83, 8).Value = 528.0714
$ws.Cells.Item(83, 9).Value = 269.14285
$ws.Cells.Item(83, 10).Value = 787
$ws.Cells.Item(83, 11).Value = 1345.71425
$ws.Cells.Item(83, 12).Value = 3935
$ws.Cells.Item(83, 13).Value = 3646.28575
$ws.Cells.Item(83, 14).Value = -13919
$ws.Cells.Item(86, 8).Value = 11765797
$ws.Cells.Item(86, 9).Value = 1175.1666
$ws.Cells.Item(86, 11).Value = 1175.1666
$ws.Cells.Item(86, 13).Value = -52.16660000000002
$ws.Cells.Item(89, 8).Value = 11765797
$ws.Cells.Item(89, 9).Value = 1175.1666
$ws.Cells.Item(89, 11).Value = 5875.833000000001
$ws.Cells.Item(89, 13).Value = -259.8330000000005
$ws.Cells.Item(99, 8).Value = 4806.0386
$ws.Cells.Item(99, 9).Value = 3845.1052
$ws.Cells.Item(99, 11).Value = 3845.1052
$ws.Cells.Item(99, 13).Value = -2347.1052
$ws.Cells.Item(107, 8).Value = 2564.0476
$ws.Cells.Item(107, 10).Value = 3924.3635
$ws.Cells.Item(107, 12).Value = 3924.3635
$ws.Cells.Item(107, 14).Value = -7764.363499999999
$ws.Cells.Item(133, 8).Value = 95000
$ws.Cells.Item(133, 10).Value = 95000
$ws.Cells.Item(133, 12).Value = 95000
$ws.Cells.Item(133, 14).Value = -105120
$ws.Cells.Item(134, 8).Value = 4179.17
$ws.Cells.Item(134, 9).Value = 4155.94
$ws.Cells.Item(134, 11).Value = 12467.82
$ws.Cells.Item(134, 13).Value = -9932.82

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 5501.1665
$ws.Cells.Item(132, 9).Value = 5501.1665
$ws.Cells.Item(132, 11).Value = 16503.4995
$ws.Cells.Item(132, 13).Value = -13973.4995
$ws.Cells.Item(141, 8).Value = 37886.5
$ws.Cells.Item(141, 10).Value = 37886.5
$ws.Cells.Item(141, 12).Value = 37886.5
$ws.Cells.Item(141, 14).Value = -48246.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 11267.111
$ws.Cells.Item(23, 9).Value = 100.833336
$ws.Cells.Item(23, 11).Value = 302.500008
$ws.Cells.Item(23, 13).Value = -67.50000799999998
$ws.Cells.Item(38, 8).Value = 481.76923
$ws.Cells.Item(38, 9).Value = 24.714285
$ws.Cells.Item(38, 10).Value = 1015
$ws.Cells.Item(38, 11).Value = 74.142855
$ws.Cells.Item(38, 12).Value = 3045
$ws.Cells.Item(38, 13).Value = 272.857145
$ws.Cells.Item(38, 14).Value = -3739
$ws.Cells.Item(114, 8).Value = 1940.4
$ws.Cells.Item(114, 9).Value = 2001.3334
$ws.Cells.Item(114, 10).Value = 1849
$ws.Cells.Item(114, 11).Value = 6004.0002
$ws.Cells.Item(114, 12).Value = 5547
$ws.Cells.Item(114, 13).Value = -2750.0002
$ws.Cells.Item(114, 14).Value = -12055
$ws.Cells.Item(129, 8).Value = 22228230
$ws.Cells.Item(129, 10).Value = 8493.25
$ws.Cells.Item(129, 12).Value = 25479.75
$ws.Cells.Item(129, 14).Value = -35479.75
$ws.Cells.Item(140, 8).Value = 4188.1816
$ws.Cells.Item(140, 9).Value = 5524
$ws.Cells.Item(140, 11).Value = 16572
$ws.Cells.Item(140, 13).Value = -11392

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(4, 8).Value = 1474.75
$ws.Cells.Item(4, 9).Value = 1474.75
$ws.Cells.Item(4, 11).Value = 1474.75
$ws.Cells.Item(4, 13).Value = -1362.75
$ws.Cells.Item(9, 8).Value = 1938.6666
$ws.Cells.Item(9, 9).Value = 1666.4
$ws.Cells.Item(9, 10).Value = 3300
$ws.Cells.Item(9, 11).Value = 1666.4
$ws.Cells.Item(9, 12).Value = 3300
$ws.Cells.Item(9, 13).Value = -1496.4
$ws.Cells.Item(9, 14).Value = -3640
$ws.Cells.Item(13, 8).Value = 1027.3334
$ws.Cells.Item(13, 9).Value = 124.333336
$ws.Cells.Item(13, 10).Value = 2833.3333
$ws.Cells.Item(13, 11).Value = 124.333336
$ws.Cells.Item(13, 12).Value = 2833.3333
$ws.Cells.Item(13, 13).Value = 14.666664
$ws.Cells.Item(13, 14).Value = -3111.3333
$ws.Cells.Item(70, 8).Value = 12724.75
$ws.Cells.Item(70, 9).Value = 2500
$ws.Cells.Item(70, 10).Value = 14185.429
$ws.Cells.Item(70, 11).Value = 2500
$ws.Cells.Item(70, 12).Value = 14185.429
$ws.Cells.Item(70, 14).Value = -14725.429
$ws.Cells.Item(70, 13).Value = -2230
$ws.Cells.Item(73, 8).Value = 12724.75
$ws.Cells.Item(73, 9).Value = 2500
$ws.Cells.Item(73, 10).Value = 14185.429
$ws.Cells.Item(73, 11).Value = 2500
$ws.Cells.Item(73, 12).Value = 14185.429
$ws.Cells.Item(73, 14).Value = -16057.429
$ws.Cells.Item(73, 13).Value = -1564
$ws.Cells.Item(97, 8).Value = 2402.6775
$ws.Cells.Item(97, 9).Value = 720.087
$ws.Cells.Item(97, 11).Value = 720.087
$ws.Cells.Item(97, 13).Value = -224.087

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3135.4348
$ws.Cells.Item(16, 9).Value = 2796.3333
$ws.Cells.Item(16, 11).Value = 2796.3333
$ws.Cells.Item(16, 13).Value = -2626.3333
$ws.Cells.Item(93, 8).Value = 9789.538
$ws.Cells.Item(93, 10).Value = 26333.777
$ws.Cells.Item(93, 12).Value = 26333.777
$ws.Cells.Item(93, 14).Value = -28829.777
$ws.Cells.Item(122, 8).Value = 5999.75
$ws.Cells.Item(122, 9).Value = 4099.6
$ws.Cells.Item(122, 10).Value = 9166.666999999999
$ws.Cells.Item(122, 11).Value = 12298.8
$ws.Cells.Item(122, 12).Value = 27500.001
$ws.Cells.Item(122, 13).Value = -9848.800000000001
$ws.Cells.Item(122, 14).Value = -32400.001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2752.6667
$ws.Cells.Item(81, 9).Value = 1442.8572
$ws.Cells.Item(81, 11).Value = 2885.7144
$ws.Cells.Item(81, 13).Value = -1824.7144
$ws.Cells.Item(84, 8).Value = 2752.6667
$ws.Cells.Item(84, 9).Value = 1442.8572
$ws.Cells.Item(84, 11).Value = 14428.572
$ws.Cells.Item(84, 13).Value = -9124.572
$ws.Cells.Item(113, 8).Value = 286.5263
$ws.Cells.Item(113, 9).Value = 298.75
$ws.Cells.Item(113, 10).Value = 221.33333
$ws.Cells.Item(113, 11).Value = 896.25
$ws.Cells.Item(113, 12).Value = 663.99999
$ws.Cells.Item(113, 13).Value = 1273.75
$ws.Cells.Item(113, 14).Value = -5003.99999
$ws.Cells.Item(122, 8).Value = 2846.9355
$ws.Cells.Item(122, 9).Value = 2187.261
$ws.Cells.Item(122, 10).Value = 4743.5
$ws.Cells.Item(122, 11).Value = 6561.782999999999
$ws.Cells.Item(122, 12).Value = 14230.5
$ws.Cells.Item(122, 13).Value = -4111.782999999999
$ws.Cells.Item(122, 14).Value = -19130.5
$ws.Cells.Item(132, 8).Value = 5464.613
$ws.Cells.Item(132, 9).Value = 4559.952
$ws.Cells.Item(132, 11).Value = 13679.856
$ws.Cells.Item(132, 13).Value = -11149.856

